# Ran code for averaged intensities on spiral schemes.
#
# This adds three new texture-sampling schemes ("Spiral-90deg-10rot-5space",
# "Spiral-90deg-15rot-5space", "Spiral-90deg-10rot-3space") to the results
# table, and repositions the existing "Gaussian-Quadrature" row so the block
# of rows reads:
#   NoRotation-tilt60deg block ... -> Gaussian-Quadrature, Spiral x3,
#   NoRotation-tilt60deg, Rotation-NoTilt, Rotation-60detTilt,
#   HexGrid-90degTilt5degRes, HexGrid-90degTilt22p5degRes,
#   HexGrid-60degTilt5degRes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: make room -------------------------------------------------
# Insert 4 new rows right before the "NoRotation-tilt60deg" row (row 10).
# This pushes the existing block (NoRotation-tilt60deg ... Gaussian-Quadrature,
# previously rows 10-16) down to rows 14-20, leaving rows 10-13 blank (and
# correctly formatted, since inserted rows inherit the formatting of the row
# above them).
$ws.Range("A10:M13").EntireRow.Insert()

# The row-insert only partially restores the thin-box-border formatting
# on the inner edges of the newly inserted block; make column A's format
# (bold, centered, full border) consistent across the new rows by pasting
# the format from the row above (row 9, which already has it).
$ws.Range("A9").Copy()
$ws.Range("A10:A13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Step 2: relocate Gaussian-Quadrature ------------------------------
# After the insert, the old "Gaussian-Quadrature" row (formerly row 16) now
# lives at row 20. Move it up to become the new row 10 (copies both values
# and formatting), then delete the now-duplicate row 20.
$ws.Range("A20:M20").Copy($ws.Range("A10:M10"))
$ws.Rows.Item(20).Delete()

# Fix up the index number for the relocated Gaussian-Quadrature row.
$ws.Cells.Item(10,1).Value = 8

# --- Step 3: fill in the new Spiral scheme rows (11-13) ----------------
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(11,3).Value = 1.002619745223633
$ws.Cells.Item(11,4).Value = 0.9675028199441925
$ws.Cells.Item(11,5).Value = 1.003725656336469
$ws.Cells.Item(11,6).Value = 1.002619745223633
$ws.Cells.Item(11,7).Value = 0.9806120319504397
$ws.Cells.Item(11,8).Value = 1.016057145405659
$ws.Cells.Item(11,9).Value = 1.003185227540635
$ws.Cells.Item(11,10).Value = 0.9675028199441925
$ws.Cells.Item(11,11).Value = 0.9856142381403308
$ws.Cells.Item(11,12).Value = 0.9941169916819821
$ws.Cells.Item(11,13).Value = 0.9956171044001714

$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(12,3).Value = 1.002489562498142
$ws.Cells.Item(12,4).Value = 0.9677584420567221
$ws.Cells.Item(12,5).Value = 1.003702095705464
$ws.Cells.Item(12,6).Value = 1.002489562498142
$ws.Cells.Item(12,7).Value = 0.9806943588813672
$ws.Cells.Item(12,8).Value = 1.015943179562734
$ws.Cells.Item(12,9).Value = 1.003114416157019
$ws.Cells.Item(12,10).Value = 0.9677584420567221
$ws.Cells.Item(12,11).Value = 0.9857302688810929
$ws.Cells.Item(12,12).Value = 0.9941099156896176
$ws.Cells.Item(12,13).Value = 0.9956170091435749

$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(13,3).Value = 1.002535898482737
$ws.Cells.Item(13,4).Value = 0.9675806455415883
$ws.Cells.Item(13,5).Value = 1.003713975749102
$ws.Cells.Item(13,6).Value = 1.002535898482737
$ws.Cells.Item(13,7).Value = 0.9805798796169275
$ws.Cells.Item(13,8).Value = 1.016060198078448
$ws.Cells.Item(13,9).Value = 1.003144958041618
$ws.Cells.Item(13,10).Value = 0.9675806455415883
$ws.Cells.Item(13,11).Value = 0.9856473106453454
$ws.Cells.Item(13,12).Value = 0.994091604564041
$ws.Cells.Item(13,13).Value = 0.9956025925850701

# --- Step 4: renumber the shifted rows (14-19) --------------------------
# Their A-column index and data were already correct after the insert+shift
# (they kept their original values), but their leading index number needs
# bumping by 4 to account for the 4 new rows above them.
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(19,1).Value = 17
